$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.734.28"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "3.755.24"
$ws.Range("E3").Value = "  -1.69%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "627.97"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").Value = "3.751.43"
$ws.Range("E7").Value = "  -1.75%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("E10").Value = "  -2.55%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.89"
$ws.Range("E12").Value = "  +4.39%  "

$ws.Range("E13").Value = "  -5.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.69"
$ws.Range("E14").Value = "  -3.75%  "

$ws.Range("D15").Value = "4.389.86"
$ws.Range("E15").Value = "  -1.53%  "

$ws.Range("D16").Value = "3.760.33"
$ws.Range("E16").Value = "  +3.76%  "

$ws.Range("D17").Value = "68.731.63"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.56"
$ws.Range("E18").Value = "  -2.56%  "

$ws.Range("E19").Value = "  -0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.97"
$ws.Range("E20").Value = "  -2.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.03"
$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.44"
$ws.Range("E22").Value = "  -2.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.700"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.50"
$ws.Range("E24").Value = "  -2.67%  "

$ws.Range("E25").Value = "  -6.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.05"
$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.09"
$ws.Range("E27").Value = "  -2.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D30").Value = "3.903.86"
$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("E31").Value = "  +1.50%  "

$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("E33").Value = "  -2.72%  "

$ws.Range("E34").Value = "  +19.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.28"
$ws.Range("E35").Value = "  -3.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").Value = "3.708.38"
$ws.Range("E37").Value = "  -1.45%  "

$ws.Range("E38").Value = "  -2.81%  "

$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("E40").Value = "  -6.00%  "

$ws.Range("E41").Value = "  -2.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.955"
$ws.Range("E43").Value = "  -2.62%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.34"
$ws.Range("E45").Value = "  +4.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "156.05"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.95"
$ws.Range("E47").Value = "  +2.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.08"
$ws.Range("E48").Value = "  +0.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.39"
$ws.Range("E49").Value = "  -3.39%  "

$ws.Range("E50").Value = "  -2.52%  "
